$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "6.360" or "29.950.35" that must stay
# literal text (multiple "." separators / significant trailing zeros). Mark
# every D cell we touch as Text first so Excel will not re-parse it as a number.
$dCells = @("D2", "D3", "D4", "D5", "D6", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.950.35'
$ws.Range("E2").Value = '  -1.47%  '
$ws.Range("D3").Value = '1.905.91'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '319.51'
$ws.Range("E5").Value = '  -2.12%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  -2.93%  '
$ws.Range("D8").Value = '0.4034'
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = '0.08251'
$ws.Range("E9").Value = '  -2.89%  '
$ws.Range("D10").Value = '41.97'
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("E11").Value = '  -2.24%  '
$ws.Range("D12").Value = '24.04'
$ws.Range("E12").Value = '  +2.27%  '
$ws.Range("D13").Value = '1.909.37'
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("D14").Value = '6.360'
$ws.Range("E14").Value = '  -1.82%  '
$ws.Range("D15").Value = '7.187'
$ws.Range("E15").Value = '  -2.53%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("D17").Value = '91.86'
$ws.Range("E17").Value = '  -3.43%  '
$ws.Range("D19").Value = '0.06494'
$ws.Range("E20").Value = '  -2.08%  '
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").Value = '5.933'
$ws.Range("E22").Value = '  -1.21%  '
$ws.Range("D23").Value = '29.986.55'
$ws.Range("E23").Value = '  -1.39%  '
$ws.Range("D24").Value = '11.23'
$ws.Range("E24").Value = '  -0.73%  '
$ws.Range("D25").Value = '2.197'
$ws.Range("E25").Value = '  -1.29%  '
$ws.Range("D26").Value = '22.16'
$ws.Range("E26").Value = '  +1.53%  '
$ws.Range("D27").Value = '2.126.01'
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").Value = '161.58'
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").Value = '2.269'
$ws.Range("E29").Value = '  -5.93%  '
$ws.Range("D30").Value = '128.60'
$ws.Range("E30").Value = '  -1.20%  '
$ws.Range("D31").Value = '1.122'
$ws.Range("E31").Value = '  +1.98%  '
$ws.Range("D32").Value = '0.1033'
$ws.Range("E32").Value = '  -2.61%  '
$ws.Range("D33").Value = '5.914'
$ws.Range("E33").Value = '  -2.17%  '
$ws.Range("D34").Value = '3.794'
$ws.Range("E34").Value = '  +0.97%  '
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value = '5.368'
$ws.Range("E35").Value = '  +2.15%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02431'
$ws.Range("E36").Value = '  -3.01%  '
$ws.Range("D37").Value = '0.06330'
$ws.Range("E37").Value = '  -4.10%  '
$ws.Range("D38").Value = '0.2139'
$ws.Range("E38").Value = '  -3.65%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = '0.6492'
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '1.191'
$ws.Range("E40").Value = '  -3.50%  '
$ws.Range("D41").Value = '8.624'
$ws.Range("E41").Value = '  -2.63%  '
$ws.Range("D42").Value = '11.29'
$ws.Range("E42").Value = '  -5.44%  '
$ws.Range("E43").Value = '  -3.02%  '
$ws.Range("D44").Value = '2.197'
$ws.Range("E44").Value = '  +6.26%  '
$ws.Range("D45").Value = '13.28'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("D46").Value = '0.6001'
$ws.Range("E46").Value = '  -2.53%  '
$ws.Range("D47").Value = '3.631'
$ws.Range("E47").Value = '  -2.33%  '
$ws.Range("D48").Value = '122.40'
$ws.Range("E48").Value = '  -2.30%  '
$ws.Range("E49").Value = '  -3.54%  '
$ws.Range("D50").Value = '78.27'
$ws.Range("E50").Value = '  -1.75%  '
$ws.Range("D51").Value = '1.128'
$ws.Range("E51").Value = '  -2.94%  '
